# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das linhas com total das tabelas

$wb = $excel.ActiveWorkbook

function Set-TextLabel {
    param($ws, $cellAddr, $text, $formatSourceAddr)
    # Writing a numeric-looking string straight into Range.Value gets coerced
    # to a number by this engine. To force a genuine text cell (matching the
    # other inline-string header cells such as B1/C1/D1) we stage the value
    # in a scratch cell formatted as Text, copy just the (now textual) value
    # into place, then restore the original cell formatting/style by pasting
    # the format from a sibling header cell that already carries it.
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $ws.Range($formatSourceAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $scratch.Clear()
}

# Sheets 1-3 and 5 use a simple "2050" label in E1; sheet 4 uses a range
# label "2041-2050" (mirroring its own C1/D1 "2015-2030" / "2031-2040"
# pattern). Sheet 6 has no E1 header at all.
$simpleLabelSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $simpleLabelSheets) {
    $ws = $wb.Worksheets.Item($name)
    Set-TextLabel $ws "E1" "2050" "D1"
}

$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextLabel $ws4 "E1" "2041-2050" "D1"

# Sheets 1-4 each have a "Total" row at row 13 that must be removed.
$totalRow13Sheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $totalRow13Sheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(13).Delete()
}

# Sheet 6 has a "Total" row at row 4 that must be removed.
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
